$d = $word.ActiveDocument

# 1. Replace the placeholder ID text (merging the trailing space run into the
#    main run so the paragraph ends with a single run).
$d.Content.Find.Execute("**ID__AFFARS_5319_topic_14__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5319_1405__ID**", 2)

# 2. Add a paragraph border (top/left/bottom/right, 5 twip space, no line)
#    and bump the left indent from 120 to 225 twips (6pt -> 11.25pt) on the
#    first paragraph.
$p = $d.Paragraphs(1)
$p.Range.Borders(-1).DistanceFromTop = 5
$p.Range.Borders(-1).DistanceFromLeft = 5
$p.Range.Borders(-1).DistanceFromBottom = 5
$p.Range.Borders(-1).DistanceFromRight = 5
$p.Range.ParagraphFormat.LeftIndent = 11.25
